$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Regenerated "K" (Strike#) values for rows 2-11 and 15-16 in column G
$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 3
    6  = 1
    7  = 0
    8  = 2
    9  = 2
    10 = 1
    11 = 1
    15 = 1
    16 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
